$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp text (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 19 de Agosto de 2020 a las 17:18"

# --- Update per-country statistics ---

# Row 4: Estados Unidos
$ws.Range("B4").Value = 5662248
$ws.Range("C4").Value = 6274
$ws.Range("D4").Value = 3012093
$ws.Range("E4").Value = 2474852
$ws.Range("G4").Value = 229
$ws.Range("H4").Value = 175303

# Row 6: India
$ws.Range("B6").Value = 2814157
$ws.Range("C6").Value = 47531
$ws.Range("D6").Value = 2075836
$ws.Range("E6").Value = 684620
$ws.Range("G6").Value = 687
$ws.Range("H6").Value = 53701

# Row 20: Italia
$ws.Range("B20").Value = 255278
$ws.Range("C20").Value = 642
$ws.Range("D20").Value = 204506
$ws.Range("E20").Value = 15360
$ws.Range("G20").Value = 7
$ws.Range("H20").Value = 35412

# Row 22: Alemania
$ws.Range("B22").Value = 228165
$ws.Range("C22").Value = 60
$ws.Range("E22").Value = 14960

# Row 50: Portugal
$ws.Range("B50").Value = 54701
$ws.Range("C50").Value = 253
$ws.Range("D50").Value = 40129
$ws.Range("E50").Value = 12786
$ws.Range("G50").Value = 2
$ws.Range("H50").Value = 1786

# Row 65: Moldavia
$ws.Range("B65").Value = 31415
$ws.Range("C65").Value = 626
$ws.Range("E65").Value = 8616
$ws.Range("G65").Value = 6
$ws.Range("H65").Value = 914

# Row 88: Noruega
$ws.Range("B88").Value = 10135
$ws.Range("C88").Value = 24
$ws.Range("E88").Value = 1016

# Row 166: Trinidad yTobago
$ws.Range("B166").Value = 650
$ws.Range("C166").Value = 21
$ws.Range("E166").Value = 498

# Row 173: Birmania
$ws.Range("B173").Value = 394
$ws.Range("C173").Value = 18
$ws.Range("D173").Value = 333
$ws.Range("E173").Value = 55

# --- Swap Montserrat / Islas Malvinas ordering (rows 213-214) ---
# They are tied on total cases (13), and the sort order flips so that
# "Islas Malvinas" now comes before "Montserrat".
$ws.Range("A213").Value = "Islas Malvinas"
$ws.Range("D213").Value = 13
$ws.Range("H213").Value = 0

$ws.Range("A214").Value = "Montserrat"
$ws.Range("D214").Value = 12
$ws.Range("H214").Value = 1
